$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with uniform run formatting) ---
# "Volume 30   Number  52" -> "Volume 30   Number  1"
$ws.Range("A8").Value = "Volume 30   Number  1"
# "Report Covering the Week  12/26/2022  Through  1/1/2023"
#  -> "Report Covering the Week  1/2/2023  Through  1/8/2023"
$ws.Range("C9").Value = "Report Covering the Week  1/2/2023  Through  1/8/2023"

# --- Precinct crime-stat table updates (rows 15-27) and historical table (rows 38, 43) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = '#,##0'
$ws.Range("K15").Value = -100
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L15").Value = -100
$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 71.428571428571
$ws.Range("I16").Value = 2
$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("J16").Value = 2
$ws.Range("J16").NumberFormat = '#,##0'
$ws.Range("K16").Value = 0
$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 3
$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("J17").Value = 2
$ws.Range("J17").NumberFormat = '#,##0'
$ws.Range("K17").Value = 50
$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = 0
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = 200
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 3
$ws.Range("I18").NumberFormat = '#,##0'
$ws.Range("J18").Value = 1
$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("K18").Value = 200
$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L18").Value = 200
$ws.Range("L18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -59.259259259259
$ws.Range("I19").Value = 5
$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("J19").Value = 10
$ws.Range("J19").NumberFormat = '#,##0'
$ws.Range("K19").Value = -50
$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = -44.444444444444
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 3
$ws.Range("J20").NumberFormat = '#,##0'
$ws.Range("K20").Value = -33.333333333333
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = 0
$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -26.315789473684
$ws.Range("F21").Value = 64
$ws.Range("H21").Value = -25.581395348837
$ws.Range("I21").Value = 15
$ws.Range("J21").Value = 19
$ws.Range("J21").NumberFormat = '#,##0'
$ws.Range("K21").Value = -21.052631578947
$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = -6.25
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 22.222222222222
$ws.Range("F24").Value = 129
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = 26.470588235294
$ws.Range("I24").Value = 33
$ws.Range("I24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 27
$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("K24").Value = 22.222222222222
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = 135.714285714286
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 120
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 105.263157894737
$ws.Range("I25").Value = 12
$ws.Range("J25").Value = 6
$ws.Range("K25").Value = 100
$ws.Range("L25").Value = 50
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 1
$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("K26").Value = -100
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = -100
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -33.333333333333
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 2
$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("J27").Value = 3
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("K27").Value = -33.333333333333
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = -33.333333333333
$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J38").Value = 87
$ws.Range("J43").Value = 919
